$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Statistics")

# Add the three new header cells (F1:H1)
$ws.Range("F1").Value = "Avg. Speed (Road 1, Direction 0)"
$ws.Range("G1").Value = "Avg. Speed (Road 1, Direction 1)"
$ws.Range("H1").Value = "Avg. Speed (Road 2, Direction 0)"

# Match the existing header formatting (bold, centered, bordered) used by A1:E1
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null

# Update the existing data row with the new simulation values
$ws.Range("A2").Value = "2024-09-01 16:16:28"
$ws.Range("B2").Value = 43.24172181965915
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 42.80508055616068
$ws.Range("E2").Value = 45.71323319191159

# Add the three new data cells (F2:H2) that go with the new headers
$ws.Range("F2").Value = 37.69324910085565
$ws.Range("G2").Value = 41.47129424928515
$ws.Range("H2").Value = 48.52575200008273
